$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Formula field update: worktime.bwAvgViews.0.* -> worktime.bwFields.0.* (BH2/BI2)
# (commit: "Update formula field in worktime"; also fixes the
# "formula override issue" described in the commit message by pointing the
# template at worktime.bwFields instead of the no-longer-valid bwAvgViews path)
$ws.Range("BH2").Value = '${worktime.bwFields.0.assyAvg}'
$ws.Range("BI2").Value = '${worktime.bwFields.0.packingAvg}'

# Output column position update: move the active selection/view over to the
# new right-hand columns (commit: "Update excel output column position")
$ws.Range("BI2").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 49
$win.ScrollRow = 1
